$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed new rows 1052-1081 by copying formatting from the last existing row (1051),
# then overwrite values/content for the Unit 36 vocabulary list.
$ws.Range("A1051:G1051").Copy($ws.Range("A1052:G1081"))

$ws.Cells.Item(1052, 1).Value = "U36_01"
$ws.Cells.Item(1052, 2).Value = 36
$ws.Cells.Item(1052, 3).Value = "Chuyến bay"
$ws.Cells.Item(1052, 4).Value = "Flight"
$ws.Cells.Item(1052, 5).Value = "International flights are currently delayed."
$ws.Cells.Item(1052, 6).Value = "An international flight"
$ws.Cells.Item(1052, 7).Value = "N"

$ws.Cells.Item(1053, 1).Value = "U36_02"
$ws.Cells.Item(1053, 2).Value = 36
$ws.Cells.Item(1053, 3).Value = "Trả lại"
$ws.Cells.Item(1053, 4).Value = "Return"
$ws.Cells.Item(1053, 5).Value = "Please return the book to the library"
$ws.Cells.Item(1053, 6).Value = "Return something to somewhere"
$ws.Cells.Item(1053, 7).Value = "V"

$ws.Cells.Item(1054, 1).Value = "U36_03"
$ws.Cells.Item(1054, 2).Value = 36
$ws.Cells.Item(1054, 3).Value = "Cất cánh"
$ws.Cells.Item(1054, 4).Value = "Take off"
$ws.Cells.Item(1054, 5).Value = "The plane will take off soon"
$ws.Cells.Item(1054, 6).Value = "a plane takes off"
$ws.Cells.Item(1054, 7).Value = "V"

$ws.Cells.Item(1055, 1).Value = "U36_04"
$ws.Cells.Item(1055, 2).Value = 36
$ws.Cells.Item(1055, 3).Value = "Cổng"
$ws.Cells.Item(1055, 4).Value = "Gate"
$ws.Cells.Item(1055, 5).Value = "Please open the gate for my car."
$ws.Cells.Item(1055, 6).Value = "open the gate"
$ws.Cells.Item(1055, 7).Value = "N"

$ws.Cells.Item(1056, 1).Value = "U36_05"
$ws.Cells.Item(1056, 2).Value = 36
$ws.Cells.Item(1056, 3).Value = "Nhập"
$ws.Cells.Item(1056, 4).Value = "Enter"
$ws.Cells.Item(1056, 5).Value = "Enter your password to log in."
$ws.Cells.Item(1056, 6).Value = "Enter a password"
$ws.Cells.Item(1056, 7).Value = "V"

$ws.Cells.Item(1057, 1).Value = "U36_06"
$ws.Cells.Item(1057, 2).Value = 36
$ws.Cells.Item(1057, 3).Value = "Xách tay"
$ws.Cells.Item(1057, 4).Value = "Carry-on"
$ws.Cells.Item(1057, 5).Value = "Carry-on luggage can be stored under the seat."
$ws.Cells.Item(1057, 6).Value = "carry-on luggage / hành lý xách tay"
$ws.Cells.Item(1057, 7).Value = "Adj"

$ws.Cells.Item(1058, 1).Value = "U36_07"
$ws.Cells.Item(1058, 2).Value = 36
$ws.Cells.Item(1058, 3).Value = "Hành lý"
$ws.Cells.Item(1058, 4).Value = "Luggage"
$ws.Cells.Item(1058, 5).Value = "You can only bring 2 pieces of luggage."
$ws.Cells.Item(1058, 6).Value = "A piece of luggage / món hành lý"
$ws.Cells.Item(1058, 7).Value = "N"

$ws.Cells.Item(1059, 1).Value = "U36_08"
$ws.Cells.Item(1059, 2).Value = 36
$ws.Cells.Item(1059, 3).Value = "Khoảng cách"
$ws.Cells.Item(1059, 4).Value = "Distance"
$ws.Cells.Item(1059, 5).Value = "The long distance will make me tired."
$ws.Cells.Item(1059, 6).Value = "Long distance / khoảng cách xa"
$ws.Cells.Item(1059, 7).Value = "N"

$ws.Cells.Item(1060, 1).Value = "U36_09"
$ws.Cells.Item(1060, 2).Value = 36
$ws.Cells.Item(1060, 3).Value = "Sự đến"
$ws.Cells.Item(1060, 4).Value = "Arrival"
$ws.Cells.Item(1060, 5).Value = "Upon arrival, they checked in at the reception."
$ws.Cells.Item(1060, 6).Value = "On or upon arrival / khi đến nơi"
$ws.Cells.Item(1060, 7).Value = "N"

$ws.Cells.Item(1061, 1).Value = "U36_10"
$ws.Cells.Item(1061, 2).Value = 36
$ws.Cells.Item(1061, 3).Value = "Hành khách"
$ws.Cells.Item(1061, 4).Value = "Passenger"
$ws.Cells.Item(1061, 5).Value = "All passenger trains are currently full."
$ws.Cells.Item(1061, 6).Value = "A passenger train / tàu chở khách"
$ws.Cells.Item(1061, 7).Value = "N"

$ws.Cells.Item(1062, 1).Value = "U36_11"
$ws.Cells.Item(1062, 2).Value = 36
$ws.Cells.Item(1062, 3).Value = "Cảng"
$ws.Cells.Item(1062, 4).Value = "Port"
$ws.Cells.Item(1062, 5).Value = "We examine the goods (kiểm tra hàng hóa) at the port of entry."
$ws.Cells.Item(1062, 6).Value = "port of entry / cảng nhập cảnh"
$ws.Cells.Item(1062, 7).Value = "N"

$ws.Cells.Item(1063, 1).Value = "U36_12"
$ws.Cells.Item(1063, 2).Value = 36
$ws.Cells.Item(1063, 3).Value = "Phà"
$ws.Cells.Item(1063, 4).Value = "Ferry"
$ws.Cells.Item(1063, 5).Value = "We take the ferry across the river every day."
$ws.Cells.Item(1063, 6).Value = "Take the ferry / đi phà"
$ws.Cells.Item(1063, 7).Value = "N"

$ws.Cells.Item(1064, 1).Value = "U36_13"
$ws.Cells.Item(1064, 2).Value = 36
$ws.Cells.Item(1064, 3).Value = "Tuyến đường"
$ws.Cells.Item(1064, 4).Value = "Route"
$ws.Cells.Item(1064, 5).Value = "The route to the hotel is dangerous."
$ws.Cells.Item(1064, 6).Value = "The route to somewhere"
$ws.Cells.Item(1064, 7).Value = "N"

$ws.Cells.Item(1065, 1).Value = "U36_14"
$ws.Cells.Item(1065, 2).Value = 36
$ws.Cells.Item(1065, 3).Value = "Thẳng"
$ws.Cells.Item(1065, 4).Value = "Direct"
$ws.Cells.Item(1065, 5).Value = "Flying from Hanoi to HCM city is a direct flight."
$ws.Cells.Item(1065, 6).Value = "A direct flight / chuyến bay thẳng"
$ws.Cells.Item(1065, 7).Value = "Adj"

$ws.Cells.Item(1066, 1).Value = "U36_15"
$ws.Cells.Item(1066, 2).Value = 36
$ws.Cells.Item(1066, 3).Value = "Tai nạn"
$ws.Cells.Item(1066, 4).Value = "Accident"
$ws.Cells.Item(1066, 5).Value = "That family had a tragic accidient last month."
$ws.Cells.Item(1066, 6).Value = "a tragic accident / tai nạn bi thảm"
$ws.Cells.Item(1066, 7).Value = "N"

$ws.Cells.Item(1067, 1).Value = "U36_16"
$ws.Cells.Item(1067, 2).Value = 36
$ws.Cells.Item(1067, 3).Value = "Tiếp tục"
$ws.Cells.Item(1067, 4).Value = "Continue"
$ws.Cells.Item(1067, 5).Value = "Continue reading on the next chapter"
$ws.Cells.Item(1067, 6).Value = "continue doing something"
$ws.Cells.Item(1067, 7).Value = "V"

$ws.Cells.Item(1068, 1).Value = "U36_17"
$ws.Cells.Item(1068, 2).Value = 36
$ws.Cells.Item(1068, 3).Value = "Dạo quanh"
$ws.Cells.Item(1068, 4).Value = "Get around"
$ws.Cells.Item(1068, 5).Value = "I need a car to get around town"
$ws.Cells.Item(1068, 6).Value = "get around town / dạo quanh thị trấn"
$ws.Cells.Item(1068, 7).Value = "V"

$ws.Cells.Item(1069, 1).Value = "U36_18"
$ws.Cells.Item(1069, 2).Value = 36
$ws.Cells.Item(1069, 3).Value = "Nhà ga"
$ws.Cells.Item(1069, 4).Value = "Station"
$ws.Cells.Item(1069, 5).Value = "Please take me to the train station"
$ws.Cells.Item(1069, 6).Value = "the train station / nhà ga tàu hỏa"
$ws.Cells.Item(1069, 7).Value = "N"

$ws.Cells.Item(1070, 1).Value = "U36_19"
$ws.Cells.Item(1070, 2).Value = 36
$ws.Cells.Item(1070, 3).Value = "Động cơ"
$ws.Cells.Item(1070, 4).Value = "Engine"
$ws.Cells.Item(1070, 5).Value = "We cannot move because the car engine is broken."
$ws.Cells.Item(1070, 6).Value = "a car engine"
$ws.Cells.Item(1070, 7).Value = "N"

$ws.Cells.Item(1071, 1).Value = "U36_20"
$ws.Cells.Item(1071, 2).Value = 36
$ws.Cells.Item(1071, 3).Value = "Bánh xe"
$ws.Cells.Item(1071, 4).Value = "Wheel"
$ws.Cells.Item(1071, 5).Value = "Please have a spare wheel for your car"
$ws.Cells.Item(1071, 6).Value = "a spare wheel / một cái bánh xe dự phòng"
$ws.Cells.Item(1071, 7).Value = "N"

$ws.Cells.Item(1072, 1).Value = "U36_21"
$ws.Cells.Item(1072, 2).Value = 36
$ws.Cells.Item(1072, 3).Value = "Tốc độ"
$ws.Cells.Item(1072, 4).Value = "Speed"
$ws.Cells.Item(1072, 5).Value = "The athlete ran at the speed of light."
$ws.Cells.Item(1072, 6).Value = "The speed of light / tốc độ ánh sáng"
$ws.Cells.Item(1072, 7).Value = "N"

$ws.Cells.Item(1073, 1).Value = "U36_22"
$ws.Cells.Item(1073, 2).Value = 36
$ws.Cells.Item(1073, 3).Value = "Mũ bảo hiểm"
$ws.Cells.Item(1073, 4).Value = "Helmet"
$ws.Cells.Item(1073, 5).Value = "Wear a safety helmet to prevent accidents"
$ws.Cells.Item(1073, 6).Value = "a safety helmet / mũ bảo hiểm"
$ws.Cells.Item(1073, 7).Value = "N"

$ws.Cells.Item(1074, 1).Value = "U36_23"
$ws.Cells.Item(1074, 2).Value = 36
$ws.Cells.Item(1074, 3).Value = "Băng qua"
$ws.Cells.Item(1074, 4).Value = "Cross"
$ws.Cells.Item(1074, 5).Value = "I helped an elderly woman cross the street"
$ws.Cells.Item(1074, 6).Value = "cross a street / sang đường"
$ws.Cells.Item(1074, 7).Value = "V"

$ws.Cells.Item(1075, 1).Value = "U36_24"
$ws.Cells.Item(1075, 2).Value = 36
$ws.Cells.Item(1075, 3).Value = "Nút giao thông"
$ws.Cells.Item(1075, 4).Value = "Intersection"
$ws.Cells.Item(1075, 5).Value = "There is traffic at the major intersection"
$ws.Cells.Item(1075, 6).Value = "A major intersection / Nút giao thông trọng yếu"
$ws.Cells.Item(1075, 7).Value = "N"

$ws.Cells.Item(1076, 1).Value = "U36_25"
$ws.Cells.Item(1076, 2).Value = 36
$ws.Cells.Item(1076, 3).Value = "Việc đỗ xe"
$ws.Cells.Item(1076, 4).Value = "Parking"
$ws.Cells.Item(1076, 5).Value = "She found a parking space for my car."
$ws.Cells.Item(1076, 6).Value = "a parking space / một chổ đậu xe"
$ws.Cells.Item(1076, 7).Value = "N"

$ws.Cells.Item(1077, 1).Value = "U36_26"
$ws.Cells.Item(1077, 2).Value = 36
$ws.Cells.Item(1077, 3).Value = "Đường quốc lộ"
$ws.Cells.Item(1077, 4).Value = "Highway"
$ws.Cells.Item(1077, 5).Value = "We are driving on the highway to Vinh"
$ws.Cells.Item(1077, 6).Value = "On the highway / trên đường quốc lộ"
$ws.Cells.Item(1077, 7).Value = "N"

$ws.Cells.Item(1078, 1).Value = "U36_27"
$ws.Cells.Item(1078, 2).Value = 36
$ws.Cells.Item(1078, 3).Value = "Hướng"
$ws.Cells.Item(1078, 4).Value = "Direction"
$ws.Cells.Item(1078, 5).Value = "The car is coming from the opposite direction"
$ws.Cells.Item(1078, 6).Value = "the opposite direction / hướng đối diện"
$ws.Cells.Item(1078, 7).Value = "N"

$ws.Cells.Item(1079, 1).Value = "U36_28"
$ws.Cells.Item(1079, 2).Value = 36
$ws.Cells.Item(1079, 3).Value = "Một chiều"
$ws.Cells.Item(1079, 4).Value = "One-way"
$ws.Cells.Item(1079, 5).Value = "The one-way traffic prevents us from turning around"
$ws.Cells.Item(1079, 6).Value = "one-way traffic / giao thông một chiều"
$ws.Cells.Item(1079, 7).Value = "Adj"

$ws.Cells.Item(1080, 1).Value = "U36_29"
$ws.Cells.Item(1080, 2).Value = 36
$ws.Cells.Item(1080, 3).Value = "Xe tải"
$ws.Cells.Item(1080, 4).Value = "Truck"
$ws.Cells.Item(1080, 5).Value = "Truck drivers need a special license"
$ws.Cells.Item(1080, 6).Value = "a truck driver / tài xế xe tải"
$ws.Cells.Item(1080, 7).Value = "N"

$ws.Cells.Item(1081, 1).Value = "U36_30"
$ws.Cells.Item(1081, 2).Value = 36
$ws.Cells.Item(1081, 3).Value = "Cái phanh"
$ws.Cells.Item(1081, 4).Value = "Brake"
$ws.Cells.Item(1081, 5).Value = "The brake pedal is broken"
$ws.Cells.Item(1081, 6).Value = "The brake pedal / bàn đạp phanh"
$ws.Cells.Item(1081, 7).Value = "N"

# Mirror the scroll position / selection left behind by Excel after entering the new rows.
$excel.ActiveWindow.ScrollRow = 1058
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C1082").Select()